$p = $ppt.ActivePresentation

# --- Slide 1 notes: update "Notes Placeholder 2" text ---
$s1 = $p.Slides.Item(1)
$notes1 = $s1.NotesPage
$notesShape1 = $notes1.Shapes.Item(2)
$notesShape1.TextFrame.TextRange.Text = "=== Original comments from ./sample-docs/pptx-copy.pptx===`nC’est quoi le clou(d) ?`n=== Original comments from ./sample-docs/pptx-copy.pptx===`nC’est quoi le clou(d) ?"

# --- Slide 2 notes: update "Notes Placeholder 2" text ---
$s2 = $p.Slides.Item(2)
$notes2 = $s2.NotesPage
$notesShape2 = $notes2.Shapes.Item(2)
$notesShape2.TextFrame.TextRange.Text = "=== Original comments from ./sample-docs/pptx-copy.pptx===`nIl faut construire un DC, un autre, …`nOn parle de Capital Expense`n`nLe mode de financement du cloud c’est de l’operational expense`nOn va payer à l’heure de consummation`nOn peut faire le parallèle avec la téléphonie portable : lorsqu’il ne fonctionne plus je le jette et j’en rachète un autre`n`nAvec quelque chose d’automatiser derrière`nOn va recréer le rack and stack par des instances`nOn appellee cela des instances car on va instancier un materiel`n`nQuelques soit votre business. Des applications qui suivent vos besoins`nAWS => Compute, storage and database `nNo upfront capital exchange`n`nEconomy at scale`nAdd and remove capacity dynamically`n`n====`n`nFor more on this, see: https://aws.amazon.com/what-is-cloud-computing/`n`n`n=== Original comments from ./sample-docs/pptx-copy.pptx===`nIl faut construire un DC, un autre, …`nOn parle de Capital Expense`n`nLe mode de financement du cloud c’est de l’operational expense`nOn va payer à l’heure de consummation`nOn peut faire le parallèle avec la téléphonie portable : lorsqu’il ne fonctionne plus je le jette et j’en rachète un autre`n`nAvec quelque chose d’automatiser derrière`nOn va recréer le rack and stack par des instances`nOn appellee cela des instances car on va instancier un materiel`n`nQuelques soit votre business. Des applications qui suivent vos besoins`nAWS => Compute, storage and database `nNo upfront capital exchange`n`nEconomy at scale`nAdd and remove capacity dynamically`n`n====`n`nFor more on this, see: https://aws.amazon.com/what-is-cloud-computing/`n`n"
